$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last 5 data rows (rows 27-31), keeping rows 1-26
$ws.Range("A27:I31").EntireRow.Delete() | Out-Null

# New data values for B2:I26 (columns VOL, B0, CET_100, CET_300, SP2_PERC, SP, SP2, SP3)
$data = @(
    @(437.435613, 199.69691, 0.048736, 3.752773, 50, 1, 32, 31),
    @(439.322651, 245.779199, -0.871437, 0.278363, 50, 1, 32, 31),
    @(425.836083, 273.494475, -0.226062, 1.401595, 46.875, 0, 30, 34),
    @(420.290737, 231.013876, -1.164374, 0.674082, 50, 1, 32, 31),
    @(417.268122, 283.137019, 0.51762, 6.655713, 46.875, 0, 30, 34),
    @(436.33015, 252.19026, -0.54052, 1.850417, 53.125, 0, 34, 30),
    @(439.253306, 251.489009, -1.123946, 3.56025, 53.125, 1, 34, 29),
    @(430.745201, 213.574831, 0.271011, 3.953429, 50, 1, 32, 31),
    @(410.238656, 300.040962, 0.061501, 3.410122, 50, 0, 32, 32),
    @(435.157637, 258.688421, -0.822244, 3.270273, 46.875, 0, 30, 34),
    @(442.729496, 263.469466, -1.052725, 1.894947, 46.875, 0, 30, 34),
    @(417.552969, 273.229674, -0.504094, 0.95021, 37.5, 0, 24, 40),
    @(429.975678, 262.793661, 0.049491, 1.928833, 50, 0, 32, 32),
    @(475.922217, 241.699599, -0.830651, 1.361503, 56.25, 1, 36, 27),
    @(417.586489, 295.863162, 0.127488, 3.029304, 46.875, 1, 30, 33),
    @(449.314289, 251.038639, -0.7132270000000001, 2.291935, 62.5, 1, 40, 23),
    @(433.419991, 243.900948, 0.262595, 1.564498, 56.25, 1, 36, 27),
    @(452.029989, 240.50265, -0.585126, 2.411265, 50, 1, 32, 31),
    @(442.851493, 234.287951, 0.36854, 3.699225, 46.875, 1, 30, 33),
    @(439.947083, 266.402591, 0.069399, 2.864446, 46.875, 1, 30, 33),
    @(420.795921, 248.201786, 0.050943, 3.516213, 50, 1, 32, 31),
    @(439.305314, 224.024904, 0.337218, 4.947861, 46.875, 1, 30, 33),
    @(438.820066, 240.491243, -1.325899, -0.793198, 53.125, 0, 34, 30),
    @(439.097307, 257.449744, 0.17731, 3.508086, 53.125, 1, 34, 29),
    @(446.819861, 289.442816, 0.868833, 9.09127, 46.875, 1, 30, 33)
)

$cols = @("B", "C", "D", "E", "F", "G", "H", "I")

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $cellref = $cols[$j] + $r
        $ws.Range($cellref).Value2 = $row[$j]
    }
}
